$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("yield")

$jasmineCells = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,27,28,30,32,33,34,35,36,39,40,41,42,43,44,48,49,50,51,52,53,54,55,56,57,58)
foreach ($r in $jasmineCells) {
    $ws.Cells.Item($r, 22).Value = "M. paniculata"
}

$curryCells = @(26)
foreach ($r in $curryCells) {
    $ws.Cells.Item($r, 22).Value = "B. koenigii"
}

$mixCells = @(29,31,37,38,45,46,47,59,60,61,62,63,64,65,66,67,68,69,70,71)
foreach ($r in $mixCells) {
    $ws.Cells.Item($r, 22).Value = "Mix"
}

$ws.Range("X26").Select() | Out-Null
